$d = $word.ActiveDocument

# 1) Underline the "Retrospective Write-up" title paragraph, including its
#    paragraph mark, so both the run and the paragraph mark (pPr/rPr) pick
#    up the <w:u w:val="single"/> formatting.
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.Font.Underline = 1

# 2) Move the hidden "_GoBack" bookmark from its old location (end of the
#    "Increasing the row size..." bug paragraph) to right after
#    "Meeting Log:" (collapsed, within that same paragraph).
$oldBookmark = $d.Bookmarks("_GoBack")
$oldBookmark.Delete()

$r = $d.Content
$found = $r.Find.Execute("Meeting Log:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$r.InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $r)
$newBookmark = $d.Bookmarks("_GoBack")
$newBookmark.Range.Text = ""
